# Add Dec 31 (Friday) as a holiday-styled column on the "16-End" sheet,
# matching the look of the other weekend/holiday columns (e.g. Dec 25 Sat,
# columns T:U) - shaded gray header/gap rows and "X" marks in every
# location's attendance row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Narrow columns AF (32) and AG (33) to match the other "X" columns (2.5 chars).
$ws.Range("AF1").ColumnWidth = 1.67
$ws.Range("AG1").ColumnWidth = 1.67

# Copy all the formatting (fills/borders/fonts) from the Dec 25 (Sat) columns
# T:U onto AF:AG -- these already carry the exact holiday-shaded styles we
# need for every row (header, tech/arrival row, per-location rows, gap rows).
$ws.Range("T2:U27").Copy()
$ws.Range("AF2:AG27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Mark attendance with "X" for every clinic-location row in the AF:AG pair.
$xRows = @(5, 6, 8, 9, 11, 12, 14, 15, 17, 18, 20, 21, 23, 24, 26, 27)
foreach ($r in $xRows) {
    $ws.Range("AF$r").Value2 = "X"
    $ws.Range("AG$r").Value2 = "X"
}
